$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily-score rows being appended to the log (2025-02-19 and 2025-02-20).
# Column A holds a date-like string ("YYYY-MM-DD"); force text formatting
# first so Excel stores it as-is instead of auto-converting it to a date
# serial number.
$ws.Range("A56:A61").NumberFormat = "@"

$data = @(
    @(56, "2025-02-19", "sleep",           $false, $true),
    @(57, "2025-02-19", "activity",        $true,  $true),
    @(58, "2025-02-19", "weekly_activity", $false, $false),
    @(59, "2025-02-20", "sleep",           $false, $true),
    @(60, "2025-02-20", "activity",        $false, $false),
    @(61, "2025-02-20", "weekly_activity", $false, $false)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 1).Value = $entry[1]
    $ws.Cells.Item($row, 2).Value = $entry[2]
    $ws.Cells.Item($row, 3).Value = $entry[3]
    $ws.Cells.Item($row, 4).Value = $entry[4]
}
